$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: new "finger exercises" entry for Lecture 18
$ws.Range("A23").Value = 45822

$ws.Range("B23").Value = 11
$ws.Range("C23").Value = 41
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 17
$ws.Range("F23").Value = "CS Introduction Lecture 18"

$ws.Range("E23").Select()
